# Update Data by bot, scripted by HH
# Target row is the single data record in row 2 of the active sheet
# (resubmission of the report: DATE_TYPE_CODE moves from annual "001" to
# quarterly "004" and the report date / financial figures refresh to the
# 2020-09-30 period).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 (DATE_TYPE_CODE) is a text code that happens to look numeric ("004"),
# so force the cell to Text before writing it, then drop back to the
# default "Normal" style so no extra number-format style sticks to the
# cell (matches the original file, which carries no style on J2).
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"
$ws.Range("J2").Style = "Normal"

# N2 (REPORT_DATE) is stored as plain text, not a real Excel date.
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Numeric financial figures for the 2020-09-30 report.
$ws.Range("O2").Value = 16109587969.96
$ws.Range("P2").Value = 2975292487.37
$ws.Range("Q2").Value = 5437588026.91
$ws.Range("R2").Value = 48.6097172766
$ws.Range("S2").Value = 174777805.61
$ws.Range("T2").Value = -17.6665885616
$ws.Range("U2").Value = 2487156485.96
$ws.Range("V2").Value = 9.6190580796
$ws.Range("W2").Value = 9016012386.16
$ws.Range("X2").Value = 2754281183.08
$ws.Range("Y2").Value = -3.7937710136
$ws.Range("Z2").Value = 1026344349.9
$ws.Range("AA2").Value = -27.9417794432
$ws.Range("AB2").Value = 7093575583.8
$ws.Range("AC2").Value = 11.6519997806
$ws.Range("AD2").Value = 14.5098993476
$ws.Range("AE2").Value = 16.863374294
$ws.Range("AF2").Value = 116.0440624415
$ws.Range("AG2").Value = 55.9667472748

Write-Output "update complete"
